# Updates cryptos list price (D) and 1h volume change (E) columns on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.177.95"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.833.43"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.79"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6652"
$ws.Range("E6").Value = "  -2.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9996"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07422"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2938"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07753"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "1.850.56"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.986"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6688"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.94"
$ws.Range("E15").Value = "  -4.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.099"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008359"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").Value = "29.166.02"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "2.084.03"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.32"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.154"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.54"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1412"
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.609"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.99"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.512"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.111"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.041"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.192"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05319"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.867"
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7481"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.140"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.645"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").Value = "1.272.96"
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01799"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.733"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9278"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.907"
$ws.Range("E42").Value = "  -2.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.08444"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.96"
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("D46").Value = "1.982.57"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5149"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.762"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.10"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05879"
$ws.Range("E51").Value = "  -0.89%  "
